$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G3").Value = 0.1018633540372671
$ws.Range("H2:H3").Value = 0.1018633540372671
$ws.Range("I2:I3").Value = 0.1577639751552795
$ws.Range("J2:J3").Value = 0.1366248496305207
$ws.Range("K2:K3").Value = 43.8
$ws.Range("L2:L3").Value = 0.1360248447204969
$ws.Range("M2:M3").Value = 17.2
$ws.Range("N2:N3").Value = 0.07263513513513513
$ws.Range("O2:O3").Value = 0.3926940639269407
$ws.Range("P2:P3").Value = 17.2
$ws.Range("Q2:Q3").Value = 0.07263513513513513
$ws.Range("R2:R3").Value = 0.3926940639269407
$ws.Range("U2:U3").Value = 20.6
$ws.Range("V2:V3").Value = 0.08699324324324324
$ws.Range("W2:W3").Value = 0.3020689655172414
$ws.Range("X2:X3").Value = 0.08846408071440458
$ws.Range("Y2:Y3").Value = 0.2136048848028368
$ws.Range("Z2:Z3").Value = 2.4788298691301
$ws.Range("AA2:AA3").Value = 0.3386697581295432
$ws.Range("AB2:AB3").Value = 0.08846408071440458
$ws.Range("AC2:AC3").Value = 0.2502056774151386
$ws.Range("AG2:AG3").Value = -20.6
$ws.Range("AJ2:AJ3").Value = -0.09528214616096208
$ws.Range("AK2:AK3").Value = -0.1402314499659632
$ws.Range("AL2:AL3").Value = 0.291
$ws.Range("AM2:AM3").Value = 0.291
$ws.Range("AO2:AO3").Value = 174.5704467353952
$ws.Range("AP2:AP3").Value = -0.3864915572232646
$ws.Range("AQ2:AQ3").Value = 174.5704467353952
